{"js": "// The document discusses file systems / disk recovery tools. The author\n// proofread the text and made three small wording corrections in the\n// first body paragraph:\n//   1. \"ostan\u00fa len sa\"        -> \"zostan\u00fa, len sa\"\n//   2. \"vie \u017ee m\u00f4\u017ee\"          -> \"vie, \u017ee m\u00f4\u017ee\"\n//   3. \"n\u00e1stroju na obnovu\"   -> \"n\u00e1stroja na obnovu\"\n//\n// Each corrected phrase is unique within the document, so we can safely\n// locate and replace each one using Word's search API.\n\nasync function replaceOnce(body, searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\nawait replaceOnce(body, \"ostan\u00fa len sa\", \"zostan\u00fa, len sa\");\nawait replaceOnce(body, \"vie \u017ee m\u00f4\u017ee\", \"vie, \u017ee m\u00f4\u017ee\");\nawait replaceOnce(body, \"n\u00e1stroju na obnovu\", \"n\u00e1stroja na obnovu\");\n", "ps1": "# The document discusses file systems / disk recovery tools. The author\n# proofread the text and made three small wording corrections in the\n# first body paragraph:\n#   1. \"ostan\u00fa len sa\"        -> \"zostan\u00fa, len sa\"\n#   2. \"vie \u017ee m\u00f4\u017ee\"          -> \"vie, \u017ee m\u00f4\u017ee\"\n#   3. \"n\u00e1stroju na obnovu\"   -> \"n\u00e1stroja na obnovu\"\n#\n# Each corrected phrase is unique within the document, so Find/Replace on\n# the whole document's Content range safely targets the right spot.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $searchText, $replacementText) {\n    $range = $doc.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replacementText\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$searchText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$replacementText, 2) | Out-Null\n}\n\nReplace-Text $d \"ostan\u00fa len sa\" \"zostan\u00fa, len sa\"\nReplace-Text $d \"vie \u017ee m\u00f4\u017ee\" \"vie, \u017ee m\u00f4\u017ee\"\nReplace-Text $d \"n\u00e1stroju na obnovu\" \"n\u00e1stroja na obnovu\"\n"}
